$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 159.5  # H11: 177.11111 -> 159.5
$ws.Cells.Item(11, 9).Value = 159.5  # I11: 177.11111 -> 159.5
$ws.Cells.Item(11, 11).Value = 159.5  # K11: 177.11111 -> 159.5
$ws.Cells.Item(11, 13).Value = -19.5  # M11: -37.11111 -> -19.5
$ws.Cells.Item(18, 8).Value = 3205.2778  # H18: 3376.2354 -> 3205.2778
$ws.Cells.Item(18, 9).Value = 2355.9375  # I18: 2493.0667 -> 2355.9375
$ws.Cells.Item(18, 11).Value = 2355.9375  # K18: 2493.0667 -> 2355.9375
$ws.Cells.Item(18, 13).Value = -2071.9375  # M18: -2209.0667 -> -2071.9375
$ws.Cells.Item(62, 8).Value = 39032.266  # H62: 39036 -> 39032.266
$ws.Cells.Item(62, 9).Value = 45390.6  # I62: 45395.08 -> 45390.6
$ws.Cells.Item(62, 11).Value = 45390.6  # K62: 45395.08 -> 45390.6
$ws.Cells.Item(62, 13).Value = -44766.6  # M62: -44771.08 -> -44766.6
$ws.Cells.Item(65, 8).Value = 39032.266  # H65: 39036 -> 39032.266
$ws.Cells.Item(65, 9).Value = 45390.6  # I65: 45395.08 -> 45390.6
$ws.Cells.Item(65, 11).Value = 226953  # K65: 226975.4 -> 226953
$ws.Cells.Item(65, 13).Value = -223833  # M65: -223855.4 -> -223833
$ws.Cells.Item(116, 8).Value = 3807.1724  # H116: 3839.75 -> 3807.1724
$ws.Cells.Item(116, 9).Value = 3699.7727  # I116: 3738.0952 -> 3699.7727
$ws.Cells.Item(116, 11).Value = 3699.7727  # K116: 3738.0952 -> 3699.7727
$ws.Cells.Item(116, 13).Value = -257.7727  # M116: -296.0952000000002 -> -257.7727
$ws.Cells.Item(132, 8).Value = 29417284  # H132: 32264066 -> 29417284
$ws.Cells.Item(132, 9).Value = 33338068  # I132: 37042236 -> 33338068
$ws.Cells.Item(132, 11).Value = 100014204  # K132: 111126708 -> 100014204
$ws.Cells.Item(132, 13).Value = -100011674  # M132: -111124178 -> -100011674
$ws.Cells.Item(135, 8).Value = 10993.272  # H135: 11158.444 -> 10993.272
$ws.Cells.Item(135, 9).Value = 7815.75  # I135: 8575.143 -> 7815.75
$ws.Cells.Item(135, 10).Value = 19466.666  # J135: 20200 -> 19466.666
$ws.Cells.Item(135, 11).Value = 70341.75  # K135: 77176.287 -> 70341.75
$ws.Cells.Item(135, 12).Value = 175199.994  # L135: 181800 -> 175199.994
$ws.Cells.Item(135, 13).Value = -67806.75  # M135: -74641.287 -> -67806.75
$ws.Cells.Item(135, 14).Value = -180269.994  # N135: -186870 -> -180269.994
$ws.Cells.Item(138, 8).Value = 503944.5  # H138: 530399.5 -> 503944.5
$ws.Cells.Item(138, 10).Value = 670842.9  # J138: 718667.4 -> 670842.9
$ws.Cells.Item(138, 12).Value = 2012528.7  # L138: 2156002.2 -> 2012528.7
$ws.Cells.Item(138, 14).Value = -2022808.7  # N138: -2166282.2 -> -2022808.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1011.7976  # H32: 1001.0941 -> 1011.7976
$ws.Cells.Item(32, 9).Value = 957.7945  # I32: 946.22974 -> 957.7945
$ws.Cells.Item(32, 11).Value = 957.7945  # K32: 946.22974 -> 957.7945
$ws.Cells.Item(32, 13).Value = -670.7945  # M32: -659.22974 -> -670.7945
$ws.Cells.Item(45, 8).Value = 4249.104  # H45: 4047.353 -> 4249.104
$ws.Cells.Item(45, 9).Value = 3901.1162  # I45: 3700.1304 -> 3901.1162
$ws.Cells.Item(45, 11).Value = 3901.1162  # K45: 3700.1304 -> 3901.1162
$ws.Cells.Item(45, 13).Value = -3524.1162  # M45: -3323.1304 -> -3524.1162
$ws.Cells.Item(74, 8).Value = 4505.2104  # H74: 4674.9443 -> 4505.2104
$ws.Cells.Item(74, 9).Value = 4899.5  # I74: 5589.4 -> 4899.5
$ws.Cells.Item(74, 11).Value = 4899.5  # K74: 5589.4 -> 4899.5
$ws.Cells.Item(74, 13).Value = -4025.5  # M74: -4715.4 -> -4025.5
$ws.Cells.Item(77, 8).Value = 4505.2104  # H77: 4674.9443 -> 4505.2104
$ws.Cells.Item(77, 9).Value = 4899.5  # I77: 5589.4 -> 4899.5
$ws.Cells.Item(77, 11).Value = 24497.5  # K77: 27947 -> 24497.5
$ws.Cells.Item(77, 13).Value = -20129.5  # M77: -23579 -> -20129.5
$ws.Cells.Item(97, 8).Value = 4210.4443  # H97: 4469.7646 -> 4210.4443
$ws.Cells.Item(97, 9).Value = 2249.4546  # I97: 2494.2 -> 2249.4546
$ws.Cells.Item(97, 11).Value = 2249.4546  # K97: 2494.2 -> 2249.4546
$ws.Cells.Item(97, 13).Value = -1753.4546  # M97: -1998.2 -> -1753.4546
$ws.Cells.Item(132, 8).Value = 4103.4067  # H132: 4254.6787 -> 4103.4067
$ws.Cells.Item(132, 9).Value = 2595.0715  # I132: 2725.4614 -> 2595.0715
$ws.Cells.Item(132, 10).Value = 4572.6665  # J132: 4717 -> 4572.6665
$ws.Cells.Item(132, 11).Value = 7785.2145  # K132: 8176.3842 -> 7785.2145
$ws.Cells.Item(132, 12).Value = 13717.9995  # L132: 14151 -> 13717.9995
$ws.Cells.Item(132, 13).Value = -5255.2145  # M132: -5646.3842 -> -5255.2145
$ws.Cells.Item(132, 14).Value = -18777.9995  # N132: -19211 -> -18777.9995
$ws.Cells.Item(135, 8).Value = 78799.8  # H135: 76666.336 -> 78799.8
$ws.Cells.Item(135, 10).Value = 78799.8  # J135: 76666.336 -> 78799.8
$ws.Cells.Item(135, 12).Value = 78799.8  # L135: 76666.336 -> 78799.8
$ws.Cells.Item(135, 14).Value = -88939.8  # N135: -86806.336 -> -88939.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1267  # H99: 1398 -> 1267
$ws.Cells.Item(99, 9).Value = 1189.3334  # I99: 1357.2 -> 1189.3334
$ws.Cells.Item(99, 11).Value = 1189.3334  # K99: 1357.2 -> 1189.3334
$ws.Cells.Item(99, 13).Value = 308.6666  # M99: 140.8 -> 308.6666
$ws.Cells.Item(105, 8).Value = 4018.8333  # H105: 4159.407 -> 4018.8333
$ws.Cells.Item(105, 9).Value = 2916.9333  # I105: 2957.75 -> 2916.9333
$ws.Cells.Item(105, 11).Value = 2916.9333  # K105: 2957.75 -> 2916.9333
$ws.Cells.Item(105, 13).Value = -1169.9333  # M105: -1210.75 -> -1169.9333
$ws.Cells.Item(134, 8).Value = 4199.5483  # H134: 4669.5186 -> 4199.5483
$ws.Cells.Item(134, 9).Value = 3121.1428  # I134: 3958.7 -> 3121.1428
$ws.Cells.Item(134, 11).Value = 9363.428400000001  # K134: 11876.1 -> 9363.428400000001
$ws.Cells.Item(134, 13).Value = -6828.428400000001  # M134: -9341.099999999999 -> -6828.428400000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2848  # H31: 2655.4 -> 2848
$ws.Cells.Item(31, 9).Value = 1175.3143  # I31: 1197.7941 -> 1175.3143
$ws.Cells.Item(31, 10).Value = 3762.75  # J31: 3406.2878 -> 3762.75
$ws.Cells.Item(31, 11).Value = 1175.3143  # K31: 1197.7941 -> 1175.3143
$ws.Cells.Item(31, 12).Value = 3762.75  # L31: 3406.2878 -> 3762.75
$ws.Cells.Item(31, 13).Value = -880.3143  # M31: -902.7941000000001 -> -880.3143
$ws.Cells.Item(31, 14).Value = -4352.75  # N31: -3996.2878 -> -4352.75
$ws.Cells.Item(34, 8).Value = 2848  # H34: 2655.4 -> 2848
$ws.Cells.Item(34, 9).Value = 1175.3143  # I34: 1197.7941 -> 1175.3143
$ws.Cells.Item(34, 10).Value = 3762.75  # J34: 3406.2878 -> 3762.75
$ws.Cells.Item(34, 11).Value = 1175.3143  # K34: 1197.7941 -> 1175.3143
$ws.Cells.Item(34, 12).Value = 3762.75  # L34: 3406.2878 -> 3762.75
$ws.Cells.Item(34, 13).Value = -973.3143  # M34: -995.7941000000001 -> -973.3143
$ws.Cells.Item(34, 14).Value = -4166.75  # N34: -3810.2878 -> -4166.75
$ws.Cells.Item(58, 8).Value = 7951.091  # H58: 7065.846 -> 7951.091
$ws.Cells.Item(58, 9).Value = 9590  # I58: 4661.3335 -> 9590
$ws.Cells.Item(58, 11).Value = 9590  # K58: 4661.3335 -> 9590
$ws.Cells.Item(58, 13).Value = -9387  # M58: -4458.3335 -> -9387
$ws.Cells.Item(132, 8).Value = 3191.8125  # H132: 3481 -> 3191.8125
$ws.Cells.Item(132, 9).Value = 2756.3635  # I132: 2940.375 -> 2756.3635
$ws.Cells.Item(132, 10).Value = 4149.8  # J132: 4562.25 -> 4149.8
$ws.Cells.Item(132, 11).Value = 8269.0905  # K132: 8821.125 -> 8269.0905
$ws.Cells.Item(132, 12).Value = 12449.4  # L132: 13686.75 -> 12449.4
$ws.Cells.Item(132, 13).Value = -5739.0905  # M132: -6291.125 -> -5739.0905
$ws.Cells.Item(132, 14).Value = -17509.4  # N132: -18746.75 -> -17509.4
$ws.Cells.Item(134, 8).Value = 3563.239  # H134: 3698.2444 -> 3563.239
$ws.Cells.Item(134, 9).Value = 3439.4187  # I134: 3555.9768 -> 3439.4187
$ws.Cells.Item(134, 10).Value = 5338  # J134: 6757 -> 5338
$ws.Cells.Item(134, 11).Value = 10318.2561  # K134: 10667.9304 -> 10318.2561
$ws.Cells.Item(134, 12).Value = 16014  # L134: 20271 -> 16014
$ws.Cells.Item(134, 13).Value = -7783.256100000001  # M134: -8132.930399999999 -> -7783.256100000001
$ws.Cells.Item(134, 14).Value = -21084  # N134: -25341 -> -21084
$ws.Cells.Item(136, 8).Value = 7951.091  # H136: 7065.846 -> 7951.091
$ws.Cells.Item(136, 9).Value = 9590  # I136: 4661.3335 -> 9590
$ws.Cells.Item(136, 11).Value = 28770  # K136: 13984.0005 -> 28770
$ws.Cells.Item(136, 13).Value = -26220  # M136: -11434.0005 -> -26220

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2714.9614  # H68: 2730.913 -> 2714.9614
$ws.Cells.Item(68, 9).Value = 1925  # I68: 1914.2858 -> 1925
$ws.Cells.Item(68, 10).Value = 3066.0557  # J68: 3088.1875 -> 3066.0557
$ws.Cells.Item(68, 11).Value = 5775  # K68: 5742.857400000001 -> 5775
$ws.Cells.Item(68, 12).Value = 9198.167099999999  # L68: 9264.5625 -> 9198.167099999999
$ws.Cells.Item(68, 13).Value = -4964  # M68: -4931.857400000001 -> -4964
$ws.Cells.Item(68, 14).Value = -10820.1671  # N68: -10886.5625 -> -10820.1671
$ws.Cells.Item(71, 8).Value = 2714.9614  # H71: 2730.913 -> 2714.9614
$ws.Cells.Item(71, 9).Value = 1925  # I71: 1914.2858 -> 1925
$ws.Cells.Item(71, 10).Value = 3066.0557  # J71: 3088.1875 -> 3066.0557
$ws.Cells.Item(71, 11).Value = 17325  # K71: 17228.5722 -> 17325
$ws.Cells.Item(71, 12).Value = 27594.5013  # L71: 27793.6875 -> 27594.5013
$ws.Cells.Item(71, 13).Value = -13269  # M71: -13172.5722 -> -13269
$ws.Cells.Item(71, 14).Value = -35706.5013  # N71: -35905.6875 -> -35706.5013
$ws.Cells.Item(107, 8).Value = 831.125  # H107: 838.7778 -> 831.125
$ws.Cells.Item(107, 10).Value = 1205.3334  # J107: 1149.8182 -> 1205.3334
$ws.Cells.Item(107, 12).Value = 3616.0002  # L107: 3449.4546 -> 3616.0002
$ws.Cells.Item(107, 14).Value = -7456.0002  # N107: -7289.4546 -> -7456.0002
$ws.Cells.Item(132, 8).Value = 1910.9166  # H132: 1602.7142 -> 1910.9166
$ws.Cells.Item(132, 9).Value = 499  # I132: 500.5 -> 499
$ws.Cells.Item(132, 10).Value = 2039.2727  # J132: 2043.6 -> 2039.2727
$ws.Cells.Item(132, 11).Value = 4491  # K132: 4504.5 -> 4491
$ws.Cells.Item(132, 12).Value = 18353.4543  # L132: 18392.4 -> 18353.4543
$ws.Cells.Item(132, 13).Value = -1961  # M132: -1974.5 -> -1961
$ws.Cells.Item(132, 14).Value = -23413.4543  # N132: -23452.4 -> -23413.4543

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6906.05  # H70: 7059.0527 -> 6906.05
$ws.Cells.Item(70, 10).Value = 7876.7856  # J70: 8175.077 -> 7876.7856
$ws.Cells.Item(70, 12).Value = 7876.7856  # L70: 8175.077 -> 7876.7856
$ws.Cells.Item(70, 14).Value = -8416.785599999999  # N70: -8715.077000000001 -> -8416.785599999999
$ws.Cells.Item(73, 8).Value = 6906.05  # H73: 7059.0527 -> 6906.05
$ws.Cells.Item(73, 10).Value = 7876.7856  # J73: 8175.077 -> 7876.7856
$ws.Cells.Item(73, 12).Value = 7876.7856  # L73: 8175.077 -> 7876.7856
$ws.Cells.Item(73, 14).Value = -9748.785599999999  # N73: -10047.077 -> -9748.785599999999
$ws.Cells.Item(102, 8).Value = 15080.825  # H102: 15855.815 -> 15080.825
$ws.Cells.Item(102, 9).Value = 1723.1852  # I102: 1832.56 -> 1723.1852
$ws.Cells.Item(102, 11).Value = 1723.1852  # K102: 1832.56 -> 1723.1852
$ws.Cells.Item(102, 13).Value = -101.1851999999999  # M102: -210.5599999999999 -> -101.1851999999999
$ws.Cells.Item(113, 8).Value = 775  # H113: 800 -> 775
$ws.Cells.Item(113, 9).Value = 775  # I113: 800 -> 775
$ws.Cells.Item(113, 11).Value = 775  # K113: 800 -> 775
$ws.Cells.Item(113, 13).Value = 1395  # M113: 1370 -> 1395
$ws.Cells.Item(132, 8).Value = 4829.706  # H132: 5262.9033 -> 4829.706
$ws.Cells.Item(132, 9).Value = 7747  # I132: 8840.846 -> 7747
$ws.Cells.Item(132, 10).Value = 2526.5789  # J132: 2678.8333 -> 2526.5789
$ws.Cells.Item(132, 11).Value = 23241  # K132: 26522.538 -> 23241
$ws.Cells.Item(132, 12).Value = 7579.736699999999  # L132: 8036.499899999999 -> 7579.736699999999
$ws.Cells.Item(132, 13).Value = -20711  # M132: -23992.538 -> -20711
$ws.Cells.Item(132, 14).Value = -12639.7367  # N132: -13096.4999 -> -12639.7367

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1938.2  # H16: 1942.92 -> 1938.2
$ws.Cells.Item(16, 9).Value = 809.3333  # I16: 817.2 -> 809.3333
$ws.Cells.Item(16, 11).Value = 809.3333  # K16: 817.2 -> 809.3333
$ws.Cells.Item(16, 13).Value = -639.3333  # M16: -647.2 -> -639.3333
$ws.Cells.Item(61, 8).Value = 22782.281  # H61: 24194.467 -> 22782.281
$ws.Cells.Item(61, 9).Value = 25519.309  # I61: 26522.12 -> 25519.309
$ws.Cells.Item(61, 10).Value = 10921.833  # J61: 12556.2 -> 10921.833
$ws.Cells.Item(61, 11).Value = 25519.309  # K61: 26522.12 -> 25519.309
$ws.Cells.Item(61, 12).Value = 10921.833  # L61: 12556.2 -> 10921.833
$ws.Cells.Item(61, 13).Value = -25317.309  # M61: -26320.12 -> -25317.309
$ws.Cells.Item(61, 14).Value = -11325.833  # N61: -12960.2 -> -11325.833
$ws.Cells.Item(113, 8).Value = 22782.281  # H113: 24194.467 -> 22782.281
$ws.Cells.Item(113, 9).Value = 25519.309  # I113: 26522.12 -> 25519.309
$ws.Cells.Item(113, 10).Value = 10921.833  # J113: 12556.2 -> 10921.833
$ws.Cells.Item(113, 11).Value = 25519.309  # K113: 26522.12 -> 25519.309
$ws.Cells.Item(113, 12).Value = 10921.833  # L113: 12556.2 -> 10921.833
$ws.Cells.Item(113, 13).Value = -23349.309  # M113: -24352.12 -> -23349.309
$ws.Cells.Item(113, 14).Value = -15261.833  # N113: -16896.2 -> -15261.833
$ws.Cells.Item(132, 8).Value = 9292.333000000001  # H132: 9115.936 -> 9292.333000000001
$ws.Cells.Item(132, 9).Value = 9308.951999999999  # I132: 9087.559999999999 -> 9308.951999999999
$ws.Cells.Item(132, 11).Value = 27926.856  # K132: 27262.68 -> 27926.856
$ws.Cells.Item(132, 13).Value = -25396.856  # M132: -24732.68 -> -25396.856
$ws.Cells.Item(133, 8).Value = 149995  # H133: 149996.5 -> 149995
$ws.Cells.Item(133, 10).Value = 149995  # J133: 149996.5 -> 149995
$ws.Cells.Item(133, 12).Value = 149995  # L133: 149996.5 -> 149995
$ws.Cells.Item(133, 14).Value = -155055  # N133: -155056.5 -> -155055
$ws.Cells.Item(136, 8).Value = 4074.6758  # H136: 4074.7026 -> 4074.6758
$ws.Cells.Item(136, 9).Value = 3893.5518  # I136: 3893.5862 -> 3893.5518
$ws.Cells.Item(136, 11).Value = 11680.6554  # K136: 11680.7586 -> 11680.6554
$ws.Cells.Item(136, 13).Value = -9130.6554  # M136: -9130.758600000001 -> -9130.6554

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 8083.6  # H126: 5669.8184 -> 8083.6
$ws.Cells.Item(126, 9).Value = 8083.6  # I126: 9467 -> 8083.6
$ws.Cells.Item(126, 10).Value = 0  # J126: 3500 -> 0
$ws.Cells.Item(126, 11).Value = 24250.8  # K126: 28401 -> 24250.8
$ws.Cells.Item(126, 12).Value = 0  # L126: 10500 -> 0
$ws.Cells.Item(126, 13).Value = -21780.8  # M126: -25931 -> -21780.8
$ws.Cells.Item(126, 14).ClearContents()  # N126 was -15440
$ws.Cells.Item(132, 8).Value = 2983.7856  # H132: 3008.6428 -> 2983.7856
$ws.Cells.Item(132, 9).Value = 2271.0588  # I132: 2312 -> 2271.0588
$ws.Cells.Item(132, 11).Value = 6813.176399999999  # K132: 6936 -> 6813.176399999999
$ws.Cells.Item(132, 13).Value = -4283.176399999999  # M132: -4406 -> -4283.176399999999
$ws.Cells.Item(136, 8).Value = 3857.5789  # H136: 4090 -> 3857.5789
$ws.Cells.Item(136, 9).Value = 2063.56  # I136: 2188.682 -> 2063.56
$ws.Cells.Item(136, 11).Value = 6190.68  # K136: 6566.045999999999 -> 6190.68
$ws.Cells.Item(136, 13).Value = -3640.68  # M136: -4016.045999999999 -> -3640.68
$ws.Cells.Item(140, 8).Value = 96663  # H140: 99995 -> 96663
$ws.Cells.Item(140, 10).Value = 96663  # J140: 99995 -> 96663
$ws.Cells.Item(140, 12).Value = 96663  # L140: 99995 -> 96663
$ws.Cells.Item(140, 14).Value = -107023  # N140: -110355 -> -107023
